# proses pembuatan laporan export
# Remove the "BAYAR SELISIH KURANG" (payment) rows from the report, keeping
# only the original "SELISIH KURANG SETORAN" (deficit) rows. The table
# shrinks from 22 data rows (B1:F22) down to 13 data rows (B1:F13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set (rows 1-4 are unchanged; rows 5-13 take over the content
# that used to live further down the old table; rows 14-22 go away).
$data = @(
    @("04/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 4/1/2023", 20000, 0),
    @("09/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 09/01/2023", 91000, 0),
    @("10/01/2023", "IQBAL", "IQBAL SELISIH KURANG SETORAN TGL 10/01/2023", 1000, 0),
    @("11/01/2023", "SANDY", "SANDY SELISIH KURANG SETORAN TGL 11/1/2023", 18000, 0),
    @("12/01/2023", "SLAMET", "SLAMET (SANDY) SELISIH KURANG SETORAN TGL 12/1/2023 ANGSURAN AN EDAH @SEWU", 25000, 0),
    @("12/01/2023", "SANDY", "SANDY YOGI SELISIH KURANG SETORAN TGL 12/01/2023", 3000, 0),
    @("18/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 18/1/2023 SALAH JUMLAH MAJELIS KIMPULAN", 1000, 0),
    @("23/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 23/01/2023 SALAH JUMLAH TABUNGAN MAJELIS SEMAR", 500, 0),
    @("24/01/2023", "YOGI", "YOGI SELISIH KURANG SETORAN TGL 24/01/2023", 1000, 0),
    @("25/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 25/01/2023", 6000, 0),
    @("26/01/2023", "IQBAL", "IQBAL SELISIH KURANG SETORAN TGL 26/01/2023", 2000, 0),
    @("30/01/2023", "ARJUNA", "ARJUNA SELISIH KURANG SETORAN TGL 30/01/2023", 38500, 0),
    @("31/01/2023", "IQBAL", "IQBAL SELISIH KURANG SETORAN ", 16000, 0)
)

# First, clear out the old extent of the table (B1:F22) so nothing stale is
# left behind once we shrink it to 13 rows.
$ws.Range("B1:F22").ClearContents()

# Column B holds dates typed in as plain text (dd/mm/yyyy); force Text
# format so Excel doesn't silently convert them to date serials.
$ws.Range("B1:B13").NumberFormat = "@"

$r = 1
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $r++
}

$ws.Range("G23").Select()
